{"js": "// Update the \"keytermCoocurrences_1964to1983\" table: each cluster row's\n// Key Terms / Size / Centrality / Density cells are replaced with the\n// redone values (per the \"Redid tables and plots\" commit).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Column indices (0-based): 0=Cluster, 1=Key Terms, 2=Size, 3=Centrality, 4=Density\n// Row indices (0-based): 0=header, 1..5=clusters 1..5\nconst updates = [\n  // row, col, newText\n  [1, 1, \"function, cardiac, human, performance, reaction time, level, skin conductance, arousal, interval, signal, anxiety\"],\n  [1, 2, \"11\"],\n  [1, 3, \"1965 (1)\"],\n  [1, 4, \"1774 (3)\"],\n\n  [2, 1, \"heart rate, feedback, control, blood pressure, cardiovascular, biofeedback, respiratory, stress, rat, technique, alpha\"],\n  [2, 2, \"11\"],\n  [2, 3, \"1860 (2)\"],\n  [2, 4, \"1743 (4)\"],\n\n  [3, 1, \"conditioning, scr, autonomic, electrodermal, habituation, stimulus, orienting response, detection, component, differential\"],\n  [3, 2, \"10\"],\n  [3, 3, \"1589 (4)\"],\n  [3, 4, \"1616 (5)\"],\n\n  [4, 1, \"sleep, eeg, recording, pattern, cognition, rem, eye movement, method, period, skin resistance\"],\n  [4, 2, \"10\"],\n  [4, 3, \"1595 (3)\"],\n  [4, 4, \"2252 (1)\"],\n\n  [5, 1, \"visual, auditory, evoked, evoked potential, behavior, child, perception, cortex\"],\n  [5, 2, \"8\"],\n  [5, 3, \"1225 (5)\"],\n  [5, 4, \"2151 (2)\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the \"keytermCoocurrences_1964to1983\" table: each cluster row's\n# Key Terms / Size / Centrality / Density cells are replaced with the\n# redone values (per the \"Redid tables and plots\" commit).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Table.Cell(row, col) is 1-indexed in the Word object model.\n# Columns: 1=Cluster, 2=Key Terms, 3=Size, 4=Centrality, 5=Density\n# Rows: 1=header, 2..6=clusters 1..5\n\n$t.Cell(2, 2).Range.Text = \"function, cardiac, human, performance, reaction time, level, skin conductance, arousal, interval, signal, anxiety\"\n$t.Cell(2, 3).Range.Text = \"11\"\n$t.Cell(2, 4).Range.Text = \"1965 (1)\"\n$t.Cell(2, 5).Range.Text = \"1774 (3)\"\n\n$t.Cell(3, 2).Range.Text = \"heart rate, feedback, control, blood pressure, cardiovascular, biofeedback, respiratory, stress, rat, technique, alpha\"\n$t.Cell(3, 3).Range.Text = \"11\"\n$t.Cell(3, 4).Range.Text = \"1860 (2)\"\n$t.Cell(3, 5).Range.Text = \"1743 (4)\"\n\n$t.Cell(4, 2).Range.Text = \"conditioning, scr, autonomic, electrodermal, habituation, stimulus, orienting response, detection, component, differential\"\n$t.Cell(4, 3).Range.Text = \"10\"\n$t.Cell(4, 4).Range.Text = \"1589 (4)\"\n$t.Cell(4, 5).Range.Text = \"1616 (5)\"\n\n$t.Cell(5, 2).Range.Text = \"sleep, eeg, recording, pattern, cognition, rem, eye movement, method, period, skin resistance\"\n$t.Cell(5, 3).Range.Text = \"10\"\n$t.Cell(5, 4).Range.Text = \"1595 (3)\"\n$t.Cell(5, 5).Range.Text = \"2252 (1)\"\n\n$t.Cell(6, 2).Range.Text = \"visual, auditory, evoked, evoked potential, behavior, child, perception, cortex\"\n$t.Cell(6, 3).Range.Text = \"8\"\n$t.Cell(6, 4).Range.Text = \"1225 (5)\"\n$t.Cell(6, 5).Range.Text = \"2151 (2)\"\n"}
